$d = $word.ActiveDocument

# Locate the end of the "MÁQUINASEscola PRO-TEC" paragraph - the text that
# must survive - and the end of the "© 2020 ... Creative Commons Attribution"
# paragraph, which (along with the blank paragraph and the "Ver no Jupiter ..."
# paragraph in between) must be removed.

$anchorRange = $d.Content
$anchorRange.Find.Execute("MÁQUINASEscola PRO-TEC", $false, $false, $false, `
    $false, $false, $true, 1, $false, "", 0) | Out-Null
$deleteStart = $anchorRange.End + 1   # step past the paragraph mark

$endRange = $d.Content
$endRange.Find.Execute( `
    "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$deleteEnd = $endRange.End + 1        # include its trailing paragraph mark

$victim = $d.Range($deleteStart, $deleteEnd)
$victim.Delete()
